# Replace the literal date string "2017-12-12" in WebParserReport!B1 with
# a real date value (serial 43070), formatted with a date number format.
# This matches the commit's intent: store a proper date instead of a
# shared-string date label (the old "Date" shared string that held the
# literal "2017-12-12" text is no longer referenced once the cell holds a
# numeric date value, so it drops out of sharedStrings.xml on save).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("WebParserReport")

$ws1.Range("B1").Value = 43070
$ws1.Range("B1").NumberFormat = "mm-dd-yy"
